$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the currency code typo: "BRA" -> "BRL" (Brazilian Real)
$ws.Range("A14").Value = "BRL"

# Update the active cell selection to A15
$ws.Range("A15").Select()
